# Rotate the codeforiati category/group columns.
#
# The sheet has columns (1-indexed): A=code, B=name, C=status,
# D=codeforiati:category-name, E=codeforiati:category-code,
# F=codeforiati:group-name, G=codeforiati:group-code
#
# The edit rotates the last four columns one step to the right for every
# row (including the header): the old G (group-code) becomes the new D,
# the old D (category-name) becomes the new E, the old E (category-code)
# becomes the new F, and the old F (group-name) becomes the new G.
#
# i.e. newD = oldG ; newE = oldD ; newF = oldE ; newG = oldF

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$used = $ws.UsedRange
$lastRow = $used.Rows.Count

for ($i = 1; $i -le $lastRow; $i++) {
    $d = $ws.Cells.Item($i, 4).Value()
    $e = $ws.Cells.Item($i, 5).Value()
    $f = $ws.Cells.Item($i, 6).Value()
    $g = $ws.Cells.Item($i, 7).Value()

    # Some of the codeforiati codes (e.g. "110", "120", ...) look like plain
    # numbers. When such values are read back through COM they come back as
    # numeric (Double) values; writing them back verbatim would turn the
    # destination cell into a genuine number instead of the original text.
    # Re-quote anything numeric with a leading apostrophe so it is written
    # back as text, exactly like the source cell.
    if ($g -is [double] -or $g -is [int]) { $g = "'" + $g.ToString() }
    if ($d -is [double] -or $d -is [int]) { $d = "'" + $d.ToString() }
    if ($e -is [double] -or $e -is [int]) { $e = "'" + $e.ToString() }
    if ($f -is [double] -or $f -is [int]) { $f = "'" + $f.ToString() }

    $ws.Cells.Item($i, 4).Value = $g
    $ws.Cells.Item($i, 5).Value = $d
    $ws.Cells.Item($i, 6).Value = $e
    $ws.Cells.Item($i, 7).Value = $f
}
